$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row => [old, new] for column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 259
$ws1.Range("F4").Value = 29
$ws1.Range("F5").Value = 740
$ws1.Range("F6").Value = 367
$ws1.Range("F10").Value = 220
$ws1.Range("F11").Value = 5958
$ws1.Range("F12").Value = 56
$ws1.Range("F13").Value = 47
$ws1.Range("F14").Value = 493
$ws1.Range("F17").Value = 357
$ws1.Range("F18").Value = 420
$ws1.Range("F21").Value = 707
$ws1.Range("F22").Value = 141
$ws1.Range("F23").Value = 97
$ws1.Range("F27").Value = 1811
$ws1.Range("F28").Value = 474

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 52
$ws2.Range("F5").Value = 270

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 233

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 233
$ws4.Range("F4").Value = 259
$ws4.Range("F5").Value = 29
$ws4.Range("F6").Value = 740
$ws4.Range("F8").Value = 367
$ws4.Range("F12").Value = 220
$ws4.Range("F13").Value = 5958
$ws4.Range("F14").Value = 56
$ws4.Range("F15").Value = 47
$ws4.Range("F17").Value = 493
$ws4.Range("F20").Value = 357
$ws4.Range("F21").Value = 420
$ws4.Range("F22").Value = 52
$ws4.Range("F25").Value = 270
$ws4.Range("F28").Value = 707
$ws4.Range("F32").Value = 141
$ws4.Range("F33").Value = 97
$ws4.Range("F37").Value = 1811
$ws4.Range("F38").Value = 474
